$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Unmerge the ranges that will be re-shaped, so values land on plain cells
# ---------------------------------------------------------------------------
$ws.Range("A3:A8").UnMerge()
$ws.Range("A9:A10").UnMerge()
$ws.Range("B3:B5").UnMerge()
$ws.Range("B6:B8").UnMerge()

# ---------------------------------------------------------------------------
# 2. Re-write the table contents (values only - existing per-cell formatting
#    for unaffected groups is preserved automatically)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "next version"
$ws.Range("B1").Value2 = "last version"
$ws.Range("C1").Value2 = "changelog"
$ws.Range("F1").Value2 = "next version"

$ws.Range("C2").Value2 = "breaking change"
$ws.Range("D2").Value2 = "feature"
$ws.Range("E2").Value2 = "bugfix"

$ws.Range("A3").Value2 = "stable"
$ws.Range("B3").Value2 = "2.0.0"
$ws.Range("C3").Value2 = "√"
$ws.Range("D3").Value2 = "—"
$ws.Range("E3").Value2 = "—"
$ws.Range("F3").Value2 = "3.0.0"

$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value2 = "×"
$ws.Range("D4").Value2 = "√"
$ws.Range("E4").Value2 = "—"
$ws.Range("F4").Value2 = "2.1.0"

$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value2 = "×"
$ws.Range("D5").Value2 = "×"
$ws.Range("E5").Value2 = "√"
$ws.Range("F5").Value2 = "2.0.1"

$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = "no stable`n version"
$ws.Range("C6").Value2 = "—"
$ws.Range("D6").Value2 = "—"
$ws.Range("E6").Value2 = "—"
$ws.Range("F6").Value2 = "1.0.0"

$ws.Range("A7").Value2 = "preview"
$ws.Range("B7").Value2 = "2.0.0"
$ws.Range("C7").Value2 = "√"
$ws.Range("D7").Value2 = "—"
$ws.Range("E7").Value2 = "—"
$ws.Range("F7").Value2 = "3.0.0b1"

$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value2 = "×"
$ws.Range("D8").Value2 = "√"
$ws.Range("E8").Value2 = "—"
$ws.Range("F8").Value2 = "2.1.0b1"

$ws.Range("A9").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value2 = "×"
$ws.Range("D9").Value2 = "×"
$ws.Range("E9").Value2 = "√"
$ws.Range("F9").Value2 = "2.0.1b1"

$ws.Range("A10").ClearContents()
$ws.Range("B10").Value2 = "2.0.0b1"
$ws.Range("C10").Value2 = "—"
$ws.Range("D10").Value2 = "—"
$ws.Range("E10").Value2 = "—"
$ws.Range("F10").Value2 = "2.0.0b2"

# ---------------------------------------------------------------------------
# 3. Re-merge into the new shape
# ---------------------------------------------------------------------------
$ws.Range("A3:A6").Merge()
$ws.Range("A7:A10").Merge()
$ws.Range("B3:B5").Merge()
$ws.Range("B7:B9").Merge()

# ---------------------------------------------------------------------------
# 4. Formatting tweaks
#    - B9/B10 used to be plain (no alignment); they now belong to the
#      centered+vertically-centered "label" family used by the other
#      merged cells in columns A/B.
#    - B6 becomes a standalone wrapped label cell.
# ---------------------------------------------------------------------------
$ws.Range("B10").HorizontalAlignment = -4108
$ws.Range("B10").VerticalAlignment = -4108

$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B6").VerticalAlignment = -4108
$ws.Range("B6").WrapText = $true
$ws.Rows(6).RowHeight = 28.8

# Column B is wider now to fit "no stable / version" / "last version"
$ws.Columns("B").ColumnWidth = 10.44

# ---------------------------------------------------------------------------
# 5. Selection, matching the author's final cursor position
# ---------------------------------------------------------------------------
$ws.Range("F13").Select()
